$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (RED, row 15) entirely - it is dropped from the data
$ws.Rows.Item(15).Delete()

# Update Quantity (column B) values for the remaining rows
$ws.Range("B2").Value = 32
$ws.Range("B3").Value = 32
$ws.Range("B4").Value = 14
$ws.Range("B5").Value = 14
$ws.Range("B6").Value = 14
$ws.Range("B7").Value = 7
$ws.Range("B8").Value = 32
$ws.Range("B9").Value = 14
$ws.Range("B10").Value = 14
$ws.Range("B11").Value = 14
$ws.Range("B12").Value = 14
$ws.Range("B13").Value = 14
$ws.Range("B14").Value = 14

# Add new column C with header "Unnamed: 2" (copy header style from B1 first)
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Unnamed: 2"

# Populate column C values
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1
